$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Range("A1").Value = "description"
$ws.Range("B1").Value = "url"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "rating"
$ws.Range("E1").Value = "page"

# E1 is a brand new column outside the original header range - copy the
# header formatting (bold font + border + alignment) from D1 onto it.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# --- Row 2: LEGO 42123 Technic McLaren Senna GTR ---
$desc2 = "LEGO 42123 Technic McLaren Senna GTR Racing Sports Collectable Model Car, Vehicle Toy Construction Set, Idea"
$url2 = "https://www.amazon.co.uk/LEGO-42123-Technic-Collectible-Construction/dp/B08G4293BD/ref=sr_1_1?crid=1SDTYPYXM70YV&keywords=lego&qid=1636379951&refinements=p_89%3Alego&rnid=1632651031&s=kids&sprefix=my+li%2Ctoys%2C163&sr=1-1"
$ws.Range("A2").Value = $desc2
$ws.Range("B2").Value = $url2
$ws.Range("C2").Value = 29.99
$ws.Range("E2").Value = 1

# --- Row 3: LEGO 76196 Marvel The Avengers Advent Calendar ---
$desc3 = "LEGO 76196 Marvel The Avengers Advent Calendar 2021 Buildable Toys with Spider-Man and Iron Man for Kids Aged 7 Idea"
$url3 = "https://www.amazon.co.uk/LEGO-76196-tbd-LSH-29-2021/dp/B08W9GQ7MV/ref=sr_1_2?crid=1SDTYPYXM70YV&keywords=lego&qid=1636379951&refinements=p_89%3Alego&rnid=1632651031&s=kids&sprefix=my+li%2Ctoys%2C163&sr=1-2"
$ws.Range("A3").Value = $desc3
$ws.Range("B3").Value = $url3
$ws.Range("C3").Value = 24.99
$ws.Range("E3").Value = 1

# --- Row 4: LEGO 41679 Friends Forest House ---
$desc4 = "LEGO 41679 Friends Forest House Toy, Treehouse Adventure Set with Mia Mini Doll and Kayak Boat Model"
$url4 = "https://www.amazon.co.uk/LEGO-41679-Friends-Treehouse-Adventure/dp/B08W5FXSQJ/ref=sr_1_3?crid=1SDTYPYXM70YV&keywords=lego&qid=1636379951&refinements=p_89%3Alego&rnid=1632651031&s=kids&sprefix=my+li%2Ctoys%2C163&sr=1-3"
$ws.Range("A4").Value = $desc4
$ws.Range("B4").Value = $url4
$ws.Range("C4").Value = 16
$ws.Range("E4").Value = 1

# --- Rating column (D2:D4): the source data stores "4.8" as text, not a
# number, so format the range as Text before typing it in (otherwise Excel
# auto-converts the numeric-looking string to a number), then drop back to
# the default cell style so no extra formatting lingers on the cells.
$ws.Range("D2:D4").NumberFormat = "@"
$ws.Range("D2").Value = "4.8"
$ws.Range("D3").Value = "4.8"
$ws.Range("D4").Value = "4.8"
$ws.Range("D2:D4").Style = "Normal"

# --- Hyperlinks on the url column (Excel auto-applies the built-in
# "Hyperlink" style - underlined, themed color - to these cells) ---
$ws.Hyperlinks.Add($ws.Range("B2"), $url2)
$ws.Hyperlinks.Add($ws.Range("B3"), $url3)
$ws.Hyperlinks.Add($ws.Range("B4"), $url4)
